$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-06 Tuesday" "2025-05-07 Wednesday"

Replace-Text "216×5=" "639×9="
Replace-Text "633×5=" "283×6="
Replace-Text "890×7=" "616×5="
Replace-Text "377×8=" "124×4="
Replace-Text "121×5=" "538×8="
Replace-Text "900×3=" "766×6="
Replace-Text "441×3=" "724×4="
Replace-Text "845×8=" "822×5="
Replace-Text "390×4=" "151×5="
Replace-Text "362×6=" "714×3="
Replace-Text "197×4=" "736×5="
Replace-Text "111×2=" "761×7="
Replace-Text "588×8=" "450×5="
Replace-Text "405×8=" "744×7="
Replace-Text "202×7=" "687×2="
Replace-Text "310×6=" "885×8="
Replace-Text "784×8=" "542×8="
Replace-Text "256×8=" "191×4="
Replace-Text "902×8=" "113×7="
Replace-Text "693×2=" "106×5="
Replace-Text "166×8=" "850×8="
Replace-Text "379×9=" "286×2="
Replace-Text "308×3=" "847×9="
Replace-Text "170×9=" "811×5="
Replace-Text "333×2=" "508×2="
